# Auto-generated Excel COM-interop edit script
# Applies cell-value updates to ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4899.6
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 3499.3333
$ws.Range("K64").Value = 7000
$ws.Range("L64").Value = 3499.3333
$ws.Range("M64").Value = -6752
$ws.Range("N64").Value = -3995.3333
$ws.Range("H67").Value = 4899.6
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 3499.3333
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 3499.3333
$ws.Range("M67").Value = -6142
$ws.Range("N67").Value = -5215.3333
$ws.Range("H74").Value = 3248584.5
$ws.Range("I74").Value = 3997296.2
$ws.Range("J74").Value = 4166.6665
$ws.Range("K74").Value = 3997296.2
$ws.Range("L74").Value = 4166.6665
$ws.Range("M74").Value = -3996360.2
$ws.Range("N74").Value = -6038.6665
$ws.Range("H76").Value = 37933896
$ws.Range("I76").Value = 40743590
$ws.Range("K76").Value = 40743590
$ws.Range("M76").Value = -40743275
$ws.Range("H77").Value = 3248584.5
$ws.Range("I77").Value = 3997296.2
$ws.Range("J77").Value = 4166.6665
$ws.Range("K77").Value = 19986481
$ws.Range("L77").Value = 20833.3325
$ws.Range("M77").Value = -19981801
$ws.Range("N77").Value = -30193.3325
$ws.Range("H79").Value = 37933896
$ws.Range("I79").Value = 40743590
$ws.Range("K79").Value = 40743590
$ws.Range("M79").Value = -40742498
$ws.Range("H116").Value = 4205.8696
$ws.Range("I116").Value = 4768.9
$ws.Range("J116").Value = 3772.7693
$ws.Range("K116").Value = 4768.9
$ws.Range("L116").Value = 3772.7693
$ws.Range("M116").Value = -1326.9
$ws.Range("N116").Value = -10656.7693
$ws.Range("H132").Value = 1745.0588
$ws.Range("I132").Value = 1856
$ws.Range("J132").Value = 598.6667
$ws.Range("K132").Value = 5568
$ws.Range("L132").Value = 1796.0001
$ws.Range("M132").Value = -3038
$ws.Range("N132").Value = -6856.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H63").Value = 1901.2941
$ws.Range("I63").Value = 1895.125
$ws.Range("K63").Value = 1895.125
$ws.Range("M63").Value = -1209.125
$ws.Range("H66").Value = 1901.2941
$ws.Range("I66").Value = 1895.125
$ws.Range("K66").Value = 9475.625
$ws.Range("M66").Value = -6043.625
$ws.Range("H74").Value = 1881.2
$ws.Range("I74").Value = 1881.2
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1881.2
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1007.2
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1881.2
$ws.Range("I77").Value = 1881.2
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9406
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -5038
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 4948.346
$ws.Range("I105").Value = 4160.8096
$ws.Range("J105").Value = 8256
$ws.Range("K105").Value = 4160.8096
$ws.Range("L105").Value = 8256
$ws.Range("M105").Value = -2413.8096
$ws.Range("N105").Value = -11750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5148.5
$ws.Range("J62").Value = 3625
$ws.Range("L62").Value = 3625
$ws.Range("N62").Value = -4873
$ws.Range("H65").Value = 5148.5
$ws.Range("J65").Value = 3625
$ws.Range("L65").Value = 18125
$ws.Range("N65").Value = -24365
$ws.Range("H86").Value = 12398
$ws.Range("I86").Value = 2996.6667
$ws.Range("K86").Value = 2996.6667
$ws.Range("M86").Value = -1873.6667
$ws.Range("H89").Value = 12398
$ws.Range("I89").Value = 2996.6667
$ws.Range("K89").Value = 14983.3335
$ws.Range("M89").Value = -9367.333500000001
$ws.Range("H134").Value = 15152804
$ws.Range("I134").Value = 1177.8518
$ws.Range("J134").Value = 83335120
$ws.Range("K134").Value = 3533.5554
$ws.Range("L134").Value = 250005360
$ws.Range("M134").Value = -998.5553999999997
$ws.Range("N134").Value = -250010430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = 48
$ws.Range("N2").Value = -426
$ws.Range("H40").Value = 9990
$ws.Range("J40").Value = 9990
$ws.Range("L40").Value = 9990
$ws.Range("N40").Value = -10292
$ws.Range("H70").Value = 9333.333000000001
$ws.Range("I70").Value = 9333.333000000001
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9333.333000000001
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9063.333000000001
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 9333.333000000001
$ws.Range("I73").Value = 9333.333000000001
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9333.333000000001
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8397.333000000001
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 3905.3914
$ws.Range("I80").Value = 3775.1333
$ws.Range("J80").Value = 4149.625
$ws.Range("K80").Value = 3775.1333
$ws.Range("L80").Value = 4149.625
$ws.Range("M80").Value = -2777.1333
$ws.Range("N80").Value = -6145.625
$ws.Range("H83").Value = 3905.3914
$ws.Range("I83").Value = 3775.1333
$ws.Range("J83").Value = 4149.625
$ws.Range("K83").Value = 18875.6665
$ws.Range("L83").Value = 20748.125
$ws.Range("M83").Value = -13883.6665
$ws.Range("N83").Value = -30732.125
$ws.Range("H113").Value = 4026.6316
$ws.Range("I113").Value = 4770.091
$ws.Range("J113").Value = 3004.375
$ws.Range("K113").Value = 4770.091
$ws.Range("L113").Value = 3004.375
$ws.Range("M113").Value = -2600.091
$ws.Range("N113").Value = -7344.375
$ws.Range("H126").Value = 4700
$ws.Range("I126").Value = 5600
$ws.Range("K126").Value = 16800
$ws.Range("M126").Value = -14330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 674569.5600000001
$ws.Range("I40").Value = 919240.25
$ws.Range("J40").Value = 1725
$ws.Range("K40").Value = 919240.25
$ws.Range("L40").Value = 1725
$ws.Range("M40").Value = -919104.25
$ws.Range("N40").Value = -1997
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3376
$ws.Range("H122").Value = 2257.1428
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -12250

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1478.5
$ws.Range("I107").Value = 1483.15
$ws.Range("J107").Value = 1455.25
$ws.Range("K107").Value = 4449.450000000001
$ws.Range("L107").Value = 4365.75
$ws.Range("M107").Value = -2529.450000000001
$ws.Range("N107").Value = -8205.75

